$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds dates stored as serial numbers. Rows 2-121
# currently hold 45177 and must be bumped to 45178.
$range = $ws.Range("C2:C121")
$range.Value = 45178
